$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5669036
$ws.Range("J17").Value = 6002397
$ws.Range("L17").Value = 18007191
$ws.Range("N17").Value = -18007527
$ws.Range("H40").Value = 1219.1333
$ws.Range("I40").Value = 758.5714
$ws.Range("K40").Value = 758.5714
$ws.Range("M40").Value = -583.5714
$ws.Range("H64").Value = 3958.25
$ws.Range("I64").Value = 3024.75
$ws.Range("J64").Value = 4425
$ws.Range("K64").Value = 3024.75
$ws.Range("L64").Value = 4425
$ws.Range("M64").Value = -2776.75
$ws.Range("N64").Value = -4921
$ws.Range("H67").Value = 3958.25
$ws.Range("I67").Value = 3024.75
$ws.Range("J67").Value = 4425
$ws.Range("K67").Value = 3024.75
$ws.Range("L67").Value = 4425
$ws.Range("M67").Value = -2166.75
$ws.Range("N67").Value = -6141
$ws.Range("H112").Value = 1071.1333
$ws.Range("J112").Value = 1083.7441
$ws.Range("L112").Value = 3251.2323
$ws.Range("N112").Value = -5467.2323
$ws.Range("H132").Value = 2796.257
$ws.Range("I132").Value = 2889.7878
$ws.Range("J132").Value = 1253
$ws.Range("K132").Value = 8669.3634
$ws.Range("L132").Value = 3759
$ws.Range("M132").Value = -6139.3634
$ws.Range("N132").Value = -8819
$ws.Range("H137").Value = 1267.081
$ws.Range("I137").Value = 1255.0769
$ws.Range("J137").Value = 1295.4546
$ws.Range("K137").Value = 3765.2307
$ws.Range("L137").Value = 3886.3638
$ws.Range("M137").Value = -1215.2307
$ws.Range("N137").Value = -8986.363799999999
$ws.Range("H141").Value = 2648.0625
$ws.Range("I141").Value = 2331.7273
$ws.Range("J141").Value = 3344
$ws.Range("K141").Value = 6995.1819
$ws.Range("L141").Value = 10032
$ws.Range("M141").Value = -1815.1819
$ws.Range("N141").Value = -20392

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2210.3333
$ws.Range("I31").Value = 2210.3333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2210.3333
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1916.3333
$ws.Range("H61").Value = 3908.7827
$ws.Range("I61").Value = 4085.4119
$ws.Range("J61").Value = 3408.3333
$ws.Range("K61").Value = 4085.4119
$ws.Range("L61").Value = 3408.3333
$ws.Range("M61").Value = -3873.4119
$ws.Range("N61").Value = -3832.3333
$ws.Range("H88").Value = 126695
$ws.Range("I88").Value = 1981.2
$ws.Range("J88").Value = 334551.34
$ws.Range("K88").Value = 1981.2
$ws.Range("L88").Value = 334551.34
$ws.Range("M88").Value = -1575.2
$ws.Range("N88").Value = -335363.34
$ws.Range("H91").Value = 126695
$ws.Range("I91").Value = 1981.2
$ws.Range("J91").Value = 334551.34
$ws.Range("K91").Value = 1981.2
$ws.Range("L91").Value = 334551.34
$ws.Range("M91").Value = -577.2
$ws.Range("N91").Value = -337359.34
$ws.Range("H132").Value = 11708.157
$ws.Range("I132").Value = 1866.7179
$ws.Range("J132").Value = 43692.832
$ws.Range("K132").Value = 5600.153700000001
$ws.Range("L132").Value = 131078.496
$ws.Range("M132").Value = -3070.153700000001
$ws.Range("N132").Value = -136138.496
$ws.Range("H136").Value = 3908.7827
$ws.Range("I136").Value = 4085.4119
$ws.Range("J136").Value = 3408.3333
$ws.Range("K136").Value = 12256.2357
$ws.Range("L136").Value = 10224.9999
$ws.Range("M136").Value = -9706.235700000001
$ws.Range("N136").Value = -15324.9999
$ws.Range("N31").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1138.6666
$ws.Range("I99").Value = 1160.3636
$ws.Range("K99").Value = 1160.3636
$ws.Range("M99").Value = 337.6364000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3677.9722
$ws.Range("I31").Value = 2869.9
$ws.Range("J31").Value = 4688.0625
$ws.Range("K31").Value = 2869.9
$ws.Range("L31").Value = 4688.0625
$ws.Range("M31").Value = -2574.9
$ws.Range("N31").Value = -5278.0625
$ws.Range("H32").Value = 2970
$ws.Range("I32").Value = 2970
$ws.Range("K32").Value = 2970
$ws.Range("M32").Value = -2654
$ws.Range("H34").Value = 3677.9722
$ws.Range("I34").Value = 2869.9
$ws.Range("J34").Value = 4688.0625
$ws.Range("K34").Value = 2869.9
$ws.Range("L34").Value = 4688.0625
$ws.Range("M34").Value = -2667.9
$ws.Range("N34").Value = -5092.0625
$ws.Range("H58").Value = 26944.95
$ws.Range("I58").Value = 1987.8
$ws.Range("J58").Value = 51902.1
$ws.Range("K58").Value = 1987.8
$ws.Range("L58").Value = 51902.1
$ws.Range("M58").Value = -1784.8
$ws.Range("N58").Value = -52308.1
$ws.Range("H134").Value = 1332.65
$ws.Range("I134").Value = 1135
$ws.Range("K134").Value = 3405
$ws.Range("M134").Value = -870
$ws.Range("H136").Value = 26944.95
$ws.Range("I136").Value = 1987.8
$ws.Range("J136").Value = 51902.1
$ws.Range("K136").Value = 5963.4
$ws.Range("L136").Value = 155706.3
$ws.Range("M136").Value = -3413.4
$ws.Range("N136").Value = -160806.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 599.75
$ws.Range("J97").Value = 999.5
$ws.Range("L97").Value = 2998.5
$ws.Range("N97").Value = -3990.5
$ws.Range("H122").Value = 604.3125
$ws.Range("I122").Value = 384
$ws.Range("J122").Value = 704.4545000000001
$ws.Range("K122").Value = 3456
$ws.Range("L122").Value = 6340.0905
$ws.Range("M122").Value = -1006
$ws.Range("N122").Value = -11240.0905
$ws.Range("H131").Value = 697.33
$ws.Range("J131").Value = 697.33
$ws.Range("L131").Value = 2091.99
$ws.Range("N131").Value = -12171.99

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15000
$ws.Range("J15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("N15").Value = -15576
$ws.Range("H81").Value = 15000
$ws.Range("J81").Value = 15000
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -16996
$ws.Range("H84").Value = 15000
$ws.Range("J84").Value = 15000
$ws.Range("L84").Value = 45000
$ws.Range("N84").Value = -54984
$ws.Range("H122").Value = 70177880
$ws.Range("I122").Value = 27779148
$ws.Range("J122").Value = 142861420
$ws.Range("K122").Value = 83337444
$ws.Range("L122").Value = 428584260
$ws.Range("M122").Value = -83334994
$ws.Range("N122").Value = -428589160
$ws.Range("H139").Value = 25614.166
$ws.Range("J139").Value = 25614.166
$ws.Range("L139").Value = 25614.166
$ws.Range("N139").Value = -35894.166

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2707.7693
$ws.Range("I22").Value = 3018.2727
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3018.2727
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -2723.2727
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 2707.7693
$ws.Range("I27").Value = 3018.2727
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 3018.2727
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2911.2727
$ws.Range("N27").Value = -1214
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H100").Value = 1764.8823
$ws.Range("I100").Value = 1264.1818
$ws.Range("J100").Value = 2682.8333
$ws.Range("K100").Value = 1264.1818
$ws.Range("L100").Value = 2682.8333
$ws.Range("M100").Value = -723.1818000000001
$ws.Range("N100").Value = -3764.8333
$ws.Range("H122").Value = 937102.9
$ws.Range("J122").Value = 3841.3635
$ws.Range("L122").Value = 11524.0905
$ws.Range("N122").Value = -16424.0905
$ws.Range("H136").Value = 1422.7667
$ws.Range("I136").Value = 1307.76
$ws.Range("J136").Value = 1997.8
$ws.Range("K136").Value = 3923.28
$ws.Range("L136").Value = 5993.4
$ws.Range("M136").Value = -1373.28
$ws.Range("N136").Value = -11093.4
$ws.Range("N38").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1470.05
$ws.Range("I132").Value = 1308.0667
$ws.Range("J132").Value = 1956
$ws.Range("K132").Value = 3924.2001
$ws.Range("L132").Value = 5868
$ws.Range("M132").Value = -1394.2001
$ws.Range("N132").Value = -10928
